# Auto-generated edit script: updates cryptos price/volume table
# to match the target snapshot (commit: "Updated cryptos list on
# Wed Mar  8 07:32:03 UTC 2023 with GitHub Actions").
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "22.007.63"
$ws.Range("E2").Value = "  -1.97%  "

$ws.Range("D3").Value = "1.553.93"
$ws.Range("E3").Value = "  -1.21%  "

$ws.Range("E4").Value = "  +0.03%  "

$ws.Range("E5").Value = "  +0.02%  "

$ws.Range("D6").Value = "'286.37"
$ws.Range("E6").Value = "  -0.59%  "

$ws.Range("D7").Value = "'0.3759"
$ws.Range("E7").Value = "  +1.73%  "

$ws.Range("D8").Value = "'0.3229"
$ws.Range("E8").Value = "  -2.78%  "

$ws.Range("D9").Value = "'1.123"
$ws.Range("E9").Value = "  -2.65%  "

$ws.Range("D10").Value = "'40.96"
$ws.Range("E10").Value = "  -13.09%  "

$ws.Range("D11").Value = "'0.07299"
$ws.Range("E11").Value = "  -2.74%  "

$ws.Range("E12").Value = "  +0.04%  "

$ws.Range("D13").Value = "'19.23"
$ws.Range("E13").Value = "  -7.42%  "

$ws.Range("D14").Value = "'5.703"
$ws.Range("E14").Value = "  -3.86%  "

$ws.Range("D15").Value = "'6.811"
$ws.Range("E15").Value = "  -1.58%  "

$ws.Range("D16").Value = "1.546.45"
$ws.Range("E16").Value = "  -0.94%  "

$ws.Range("D17").Value = "'0.00001078"
$ws.Range("E17").Value = "  -3.33%  "

$ws.Range("D19").Value = "'84.90"
$ws.Range("E19").Value = "  -3.99%  "

$ws.Range("D20").Value = "'6.417"
$ws.Range("E20").Value = "  +0.50%  "

$ws.Range("E21").Value = "  +0.07%  "

$ws.Range("D22").Value = "'15.98"
$ws.Range("E22").Value = "  -3.11%  "

$ws.Range("D23").Value = "'11.53"
$ws.Range("E23").Value = "  -3.81%  "

$ws.Range("D24").Value = "22.018.81"
$ws.Range("E24").Value = "  -1.89%  "

$ws.Range("D25").Value = "'2.233"
$ws.Range("E25").Value = "  -6.40%  "

$ws.Range("D26").Value = "'2.519"
$ws.Range("E26").Value = "  -4.33%  "

$ws.Range("D27").Value = "'149.99"
$ws.Range("E27").Value = "  -0.58%  "

$ws.Range("D28").Value = "'18.83"
$ws.Range("E28").Value = "  -3.87%  "

$ws.Range("D29").Value = "'4.845"
$ws.Range("E29").Value = "  -2.37%  "

$ws.Range("D30").Value = "1.726.26"
$ws.Range("E30").Value = "  -0.80%  "

$ws.Range("D31").Value = "'119.88"
$ws.Range("E31").Value = "  -4.10%  "

$ws.Range("D32").Value = "'1.120"
$ws.Range("E32").Value = "  +2.07%  "

$ws.Range("D33").Value = "'5.898"
$ws.Range("E33").Value = "  -3.01%  "

$ws.Range("B34").Value = "Stellar"
$ws.Range("C34").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D34").Value = "'0.08166"
$ws.Range("E34").Value = "  -1.60%  "

$ws.Range("D35").Value = "'9.256"
$ws.Range("E35").Value = "  -6.22%  "

$ws.Range("B36").Value = "WEMIXTOKEN"
$ws.Range("C36").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D36").Value = "'1.639"
$ws.Range("E36").Value = "  -17.61%  "

$ws.Range("D37").Value = "'5.223"
$ws.Range("E37").Value = "  -2.00%  "

$ws.Range("D38").Value = "'0.02280"
$ws.Range("E38").Value = "  -6.78%  "

$ws.Range("D39").Value = "'0.06145"
$ws.Range("E39").Value = "  -3.81%  "

$ws.Range("D40").Value = "'0.2109"
$ws.Range("E40").Value = "  -4.94%  "

$ws.Range("D41").Value = "'1.214"
$ws.Range("E41").Value = "  -7.38%  "

$ws.Range("D42").Value = "'10.87"
$ws.Range("E42").Value = "  -4.61%  "

$ws.Range("E43").Value = "  +0.05%  "

$ws.Range("D44").Value = "'0.5928"
$ws.Range("E44").Value = "  -5.05%  "

$ws.Range("D45").Value = "'13.51"
$ws.Range("E45").Value = "  -3.27%  "

$ws.Range("D46").Value = "'3.721"
$ws.Range("E46").Value = "  -1.38%  "

$ws.Range("D47").Value = "'0.5734"
$ws.Range("E47").Value = "  -5.37%  "

$ws.Range("D48").Value = "'1.933"
$ws.Range("E48").Value = "  -5.47%  "

$ws.Range("D49").Value = "'119.95"
$ws.Range("E49").Value = "  -3.87%  "

$ws.Range("D50").Value = "'1.152"
$ws.Range("E50").Value = "  -4.33%  "

$ws.Range("D51").Value = "'0.06920"
$ws.Range("E51").Value = "  -3.81%  "
